# TradingModel_v2 - 2021/11/17 data updated
# Previously the last row (row 10, date 2021-11-16) carried the "date only"
# number format (s=3). Now that a new row is appended for 2021-11-17, that
# distinct format moves down to the new last row, and row 10 reverts to the
# regular datetime format used by the rest of the date column (s=2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (2021-11-16) reverts to the standard date/time format.
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 11 (2021-11-17) data.
$ws.Range("A11").Value = 44517
$ws.Range("A11").NumberFormat = "YYYY-MM-DD"
$ws.Range("B11").Value = 66336.55
